# Apply the edit described by the diff:
# - On the "März" sheet, fill in a new work interval on row 27
#   (13:00 - 15:30, "Literaturrecherche und Coding") and on row 32
#   (13:15 - 15:30, "Coding"). All downstream formulas (monthly /
#   yearly totals) recompute automatically.
# - Leave the "März" sheet as the active sheet/selection (E32), matching
#   the author's last interaction before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("März")

# Times are stored as fractions of a day.
$ws.Range("D27").Value = [double]13 / 24          # 13:00
$ws.Range("E27").Value = [double](15*60+30) / 1440 # 15:30
$ws.Range("O27").Value = "Literaturrecherche und Coding"

$ws.Range("D32").Value = [double](13*60+15) / 1440 # 13:15
$ws.Range("E32").Value = [double](15*60+30) / 1440 # 15:30
$ws.Range("O32").Value = "Coding"

# Make März the active sheet, with E32 selected (matches author's
# final cursor position before save).
$ws.Activate()
$ws.Range("E32").Select()
